# Insert a new data row above row 128 (this pushes the existing rows 128-223
# down to 129-224, preserving all of their data/formatting), then populate the
# newly inserted row 128 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("128:128").Insert()

$ws.Cells.Item(128, 1).Value  = 10
$ws.Cells.Item(128, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(128, 3).Value  = "La Araucanía"
$ws.Cells.Item(128, 4).Value  = 44512
$ws.Cells.Item(128, 5).Value  = 9
$ws.Cells.Item(128, 6).Value  = 100114013
$ws.Cells.Item(128, 7).Value  = "Zanahoria"
$ws.Cells.Item(128, 8).Value  = "Sin especificar"
$ws.Cells.Item(128, 9).Value  = "Primera"
$ws.Cells.Item(128, 10).Value = 80
$ws.Cells.Item(128, 11).Value = 12000
$ws.Cells.Item(128, 12).Value = 12000
$ws.Cells.Item(128, 13).Value = 12000
$ws.Cells.Item(128, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(128, 15).Value = "Región del Maule"
$ws.Cells.Item(128, 16).Value = 600
$ws.Cells.Item(128, 17).Value = 20
$ws.Cells.Item(128, 18).Value = "Hortaliza"
